$d = $word.ActiveDocument
$count = 0
$found = $d.Content.Find.Execute("2024-05-21 Tuesday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-05-22 Wednesday", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 2024-05-21 Tuesday" }
$found = $d.Content.Find.Execute("82-59=23", $true, $true, $false, $false, $false, $true, 1, $false, "38+38=76", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 82-59=23" }
$found = $d.Content.Find.Execute("81-36=45", $true, $true, $false, $false, $false, $true, 1, $false, "19+12=31", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 81-36=45" }
$found = $d.Content.Find.Execute("72-67=5", $true, $true, $false, $false, $false, $true, 1, $false, "63-8=55", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 72-67=5" }
$found = $d.Content.Find.Execute("24+67=91", $true, $true, $false, $false, $false, $true, 1, $false, "92-55=37", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 24+67=91" }
$found = $d.Content.Find.Execute("63-45=18", $true, $true, $false, $false, $false, $true, 1, $false, "93-14=79", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 63-45=18" }
$found = $d.Content.Find.Execute("91-67=24", $true, $true, $false, $false, $false, $true, 1, $false, "16+19=35", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 91-67=24" }
$found = $d.Content.Find.Execute("86-29=57", $true, $true, $false, $false, $false, $true, 1, $false, "57+34=91", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 86-29=57" }
$found = $d.Content.Find.Execute("66+8=74", $true, $true, $false, $false, $false, $true, 1, $false, "26+45=71", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 66+8=74" }
$found = $d.Content.Find.Execute("61-36=25", $true, $true, $false, $false, $false, $true, 1, $false, "94-65=29", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 61-36=25" }
$found = $d.Content.Find.Execute("51-37=14", $true, $true, $false, $false, $false, $true, 1, $false, "52-24=28", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 51-37=14" }
$found = $d.Content.Find.Execute("37+48=85", $true, $true, $false, $false, $false, $true, 1, $false, "24-7=17", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 37+48=85" }
$found = $d.Content.Find.Execute("14+48=62", $true, $true, $false, $false, $false, $true, 1, $false, "79+9=88", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 14+48=62" }
$found = $d.Content.Find.Execute("29+34=63", $true, $true, $false, $false, $false, $true, 1, $false, "96-47=49", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 29+34=63" }
$found = $d.Content.Find.Execute("66+18=84", $true, $true, $false, $false, $false, $true, 1, $false, "9+34=43", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 66+18=84" }
$found = $d.Content.Find.Execute("25+29=54", $true, $true, $false, $false, $false, $true, 1, $false, "59+18=77", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 25+29=54" }
$found = $d.Content.Find.Execute("77+4=81", $true, $true, $false, $false, $false, $true, 1, $false, "23+68=91", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 77+4=81" }
$found = $d.Content.Find.Execute("8+65=73", $true, $true, $false, $false, $false, $true, 1, $false, "43-25=18", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 8+65=73" }
$found = $d.Content.Find.Execute("13-4=9", $true, $true, $false, $false, $false, $true, 1, $false, "35+6=41", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 13-4=9" }
$found = $d.Content.Find.Execute("80-61=19", $true, $true, $false, $false, $false, $true, 1, $false, "85-8=77", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 80-61=19" }
$found = $d.Content.Find.Execute("86-27=59", $true, $true, $false, $false, $false, $true, 1, $false, "94-85=9", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 86-27=59" }
$found = $d.Content.Find.Execute("36+28=64", $true, $true, $false, $false, $false, $true, 1, $false, "92-19=73", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 36+28=64" }
$found = $d.Content.Find.Execute("71-54=17", $true, $true, $false, $false, $false, $true, 1, $false, "25-6=19", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 71-54=17" }
$found = $d.Content.Find.Execute("94-87=7", $true, $true, $false, $false, $false, $true, 1, $false, "64+28=92", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 94-87=7" }
$found = $d.Content.Find.Execute("92-38=54", $true, $true, $false, $false, $false, $true, 1, $false, "7+26=33", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 92-38=54" }
$found = $d.Content.Find.Execute("94-28=66", $true, $true, $false, $false, $false, $true, 1, $false, "46+49=95", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 94-28=66" }
$found = $d.Content.Find.Execute("17+35=52", $true, $true, $false, $false, $false, $true, 1, $false, "32-28=4", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 17+35=52" }
$found = $d.Content.Find.Execute("75-19=56", $true, $true, $false, $false, $false, $true, 1, $false, "55+37=92", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 75-19=56" }
$found = $d.Content.Find.Execute("24+57=81", $true, $true, $false, $false, $false, $true, 1, $false, "92-55=37", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 24+57=81" }
$found = $d.Content.Find.Execute("34-18=16", $true, $true, $false, $false, $false, $true, 1, $false, "39+49=88", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 34-18=16" }
$found = $d.Content.Find.Execute("26+7=33", $true, $true, $false, $false, $false, $true, 1, $false, "3+8=11", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 26+7=33" }
$found = $d.Content.Find.Execute("58+13=71", $true, $true, $false, $false, $false, $true, 1, $false, "28+13=41", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 58+13=71" }
$found = $d.Content.Find.Execute("73-45=28", $true, $true, $false, $false, $false, $true, 1, $false, "35+16=51", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 73-45=28" }
$found = $d.Content.Find.Execute("51-49=2", $true, $true, $false, $false, $false, $true, 1, $false, "5+76=81", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 51-49=2" }
$found = $d.Content.Find.Execute("16+7=23", $true, $true, $false, $false, $false, $true, 1, $false, "26+15=41", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 16+7=23" }
$found = $d.Content.Find.Execute("15+19=34", $true, $true, $false, $false, $false, $true, 1, $false, "16+79=95", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 15+19=34" }
$found = $d.Content.Find.Execute("62+9=71", $true, $true, $false, $false, $false, $true, 1, $false, "66+5=71", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 62+9=71" }
$found = $d.Content.Find.Execute("28+35=63", $true, $true, $false, $false, $false, $true, 1, $false, "86-57=29", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 28+35=63" }
$found = $d.Content.Find.Execute("38+56=94", $true, $true, $false, $false, $false, $true, 1, $false, "31-3=28", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 38+56=94" }
$found = $d.Content.Find.Execute("44+48=92", $true, $true, $false, $false, $false, $true, 1, $false, "58+7=65", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 44+48=92" }
$found = $d.Content.Find.Execute("84-26=58", $true, $true, $false, $false, $false, $true, 1, $false, "8+76=84", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 84-26=58" }
$found = $d.Content.Find.Execute("52-23=29", $true, $true, $false, $false, $false, $true, 1, $false, "39+57=96", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 52-23=29" }
$found = $d.Content.Find.Execute("75-48=27", $true, $true, $false, $false, $false, $true, 1, $false, "94-27=67", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 75-48=27" }
$found = $d.Content.Find.Execute("24-5=19", $true, $true, $false, $false, $false, $true, 1, $false, "92-48=44", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 24-5=19" }
$found = $d.Content.Find.Execute("23-17=6", $true, $true, $false, $false, $false, $true, 1, $false, "67-29=38", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 23-17=6" }
$found = $d.Content.Find.Execute("7+8=15", $true, $true, $false, $false, $false, $true, 1, $false, "87-48=39", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 7+8=15" }
$found = $d.Content.Find.Execute("73-58=15", $true, $true, $false, $false, $false, $true, 1, $false, "32-24=8", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 73-58=15" }
$found = $d.Content.Find.Execute("71-48=23", $true, $true, $false, $false, $false, $true, 1, $false, "43+48=91", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 71-48=23" }
$found = $d.Content.Find.Execute("26+26=52", $true, $true, $false, $false, $false, $true, 1, $false, "35+28=63", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 26+26=52" }
$found = $d.Content.Find.Execute("61-39=22", $true, $true, $false, $false, $false, $true, 1, $false, "82-64=18", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 61-39=22" }
$found = $d.Content.Find.Execute("40-28=12", $true, $true, $false, $false, $false, $true, 1, $false, "86-17=69", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 40-28=12" }
$found = $d.Content.Find.Execute("71-12=59", $true, $true, $false, $false, $false, $true, 1, $false, "40-37=3", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 71-12=59" }
$found = $d.Content.Find.Execute("18+34=52", $true, $true, $false, $false, $false, $true, 1, $false, "37+38=75", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 18+34=52" }
$found = $d.Content.Find.Execute("44-8=36", $true, $true, $false, $false, $false, $true, 1, $false, "60-31=29", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 44-8=36" }
$found = $d.Content.Find.Execute("38+18=56", $true, $true, $false, $false, $false, $true, 1, $false, "48+44=92", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 38+18=56" }
$found = $d.Content.Find.Execute("71-25=46", $true, $true, $false, $false, $false, $true, 1, $false, "58+15=73", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 71-25=46" }
$found = $d.Content.Find.Execute("90-76=14", $true, $true, $false, $false, $false, $true, 1, $false, "92-64=28", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 90-76=14" }
$found = $d.Content.Find.Execute("68+13=81", $true, $true, $false, $false, $false, $true, 1, $false, "34+7=41", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 68+13=81" }
$found = $d.Content.Find.Execute("93-28=65", $true, $true, $false, $false, $false, $true, 1, $false, "91-8=83", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 93-28=65" }
$found = $d.Content.Find.Execute("42-35=7", $true, $true, $false, $false, $false, $true, 1, $false, "69+18=87", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 42-35=7" }
$found = $d.Content.Find.Execute("29+52=81", $true, $true, $false, $false, $false, $true, 1, $false, "6+17=23", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 29+52=81" }
$found = $d.Content.Find.Execute("17+6=23", $true, $true, $false, $false, $false, $true, 1, $false, "55+29=84", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 17+6=23" }
$found = $d.Content.Find.Execute("83-77=6", $true, $true, $false, $false, $false, $true, 1, $false, "36+57=93", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 83-77=6" }
$found = $d.Content.Find.Execute("37+9=46", $true, $true, $false, $false, $false, $true, 1, $false, "82-18=64", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 37+9=46" }
$found = $d.Content.Find.Execute("48+18=66", $true, $true, $false, $false, $false, $true, 1, $false, "68-39=29", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 48+18=66" }
$found = $d.Content.Find.Execute("38+25=63", $true, $true, $false, $false, $false, $true, 1, $false, "90-73=17", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 38+25=63" }
$found = $d.Content.Find.Execute("84-77=7", $true, $true, $false, $false, $false, $true, 1, $false, "82-3=79", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 84-77=7" }
$found = $d.Content.Find.Execute("16+75=91", $true, $true, $false, $false, $false, $true, 1, $false, "40-24=16", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 16+75=91" }
$found = $d.Content.Find.Execute("87-78=9", $true, $true, $false, $false, $false, $true, 1, $false, "80-21=59", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 87-78=9" }
$found = $d.Content.Find.Execute("64+8=72", $true, $true, $false, $false, $false, $true, 1, $false, "57+4=61", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 64+8=72" }
$found = $d.Content.Find.Execute("54+27=81", $true, $true, $false, $false, $false, $true, 1, $false, "46-29=17", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 54+27=81" }
$found = $d.Content.Find.Execute("23+8=31", $true, $true, $false, $false, $false, $true, 1, $false, "83-58=25", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 23+8=31" }
$found = $d.Content.Find.Execute("91-53=38", $true, $true, $false, $false, $false, $true, 1, $false, "44-39=5", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 91-53=38" }
$found = $d.Content.Find.Execute("89+9=98", $true, $true, $false, $false, $false, $true, 1, $false, "55-7=48", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 89+9=98" }
$found = $d.Content.Find.Execute("68+19=87", $true, $true, $false, $false, $false, $true, 1, $false, "77+5=82", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 68+19=87" }
$found = $d.Content.Find.Execute("7+28=35", $true, $true, $false, $false, $false, $true, 1, $false, "39+12=51", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 7+28=35" }
$found = $d.Content.Find.Execute("90-6=84", $true, $true, $false, $false, $false, $true, 1, $false, "38+33=71", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 90-6=84" }
$found = $d.Content.Find.Execute("56-28=28", $true, $true, $false, $false, $false, $true, 1, $false, "23-6=17", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 56-28=28" }
$found = $d.Content.Find.Execute("28+69=97", $true, $true, $false, $false, $false, $true, 1, $false, "55-37=18", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 28+69=97" }
$found = $d.Content.Find.Execute("65-17=48", $true, $true, $false, $false, $false, $true, 1, $false, "20-5=15", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 65-17=48" }
$found = $d.Content.Find.Execute("20-19=1", $true, $true, $false, $false, $false, $true, 1, $false, "17+6=23", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 20-19=1" }
$found = $d.Content.Find.Execute("51-2=49", $true, $true, $false, $false, $false, $true, 1, $false, "18+48=66", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 51-2=49" }
$found = $d.Content.Find.Execute("44-16=28", $true, $true, $false, $false, $false, $true, 1, $false, "9+82=91", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 44-16=28" }
$found = $d.Content.Find.Execute("70-17=53", $true, $true, $false, $false, $false, $true, 1, $false, "80-9=71", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 70-17=53" }
$found = $d.Content.Find.Execute("16+69=85", $true, $true, $false, $false, $false, $true, 1, $false, "61-44=17", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 16+69=85" }
$found = $d.Content.Find.Execute("62-54=8", $true, $true, $false, $false, $false, $true, 1, $false, "85+6=91", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 62-54=8" }
$found = $d.Content.Find.Execute("82-55=27", $true, $true, $false, $false, $false, $true, 1, $false, "33-14=19", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 82-55=27" }
$found = $d.Content.Find.Execute("52+9=61", $true, $true, $false, $false, $false, $true, 1, $false, "57+18=75", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 52+9=61" }
$found = $d.Content.Find.Execute("43+49=92", $true, $true, $false, $false, $false, $true, 1, $false, "39+33=72", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 43+49=92" }
$found = $d.Content.Find.Execute("85-9=76", $true, $true, $false, $false, $false, $true, 1, $false, "5+38=43", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 85-9=76" }
$found = $d.Content.Find.Execute("96-38=58", $true, $true, $false, $false, $false, $true, 1, $false, "56+16=72", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 96-38=58" }
$found = $d.Content.Find.Execute("16+57=73", $true, $true, $false, $false, $false, $true, 1, $false, "83-16=67", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 16+57=73" }
$found = $d.Content.Find.Execute("63-25=38", $true, $true, $false, $false, $false, $true, 1, $false, "28+13=41", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 63-25=38" }
$found = $d.Content.Find.Execute("47+29=76", $true, $true, $false, $false, $false, $true, 1, $false, "34-16=18", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 47+29=76" }
$found = $d.Content.Find.Execute("30-26=4", $true, $true, $false, $false, $false, $true, 1, $false, "30-11=19", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 30-26=4" }
$found = $d.Content.Find.Execute("74-65=9", $true, $true, $false, $false, $false, $true, 1, $false, "18+53=71", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 74-65=9" }
$found = $d.Content.Find.Execute("34-15=19", $true, $true, $false, $false, $false, $true, 1, $false, "45+47=92", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 34-15=19" }
$found = $d.Content.Find.Execute("25+17=42", $true, $true, $false, $false, $false, $true, 1, $false, "70-51=19", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 25+17=42" }
$found = $d.Content.Find.Execute("2+89=91", $true, $true, $false, $false, $false, $true, 1, $false, "28+45=73", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 2+89=91" }
$found = $d.Content.Find.Execute("63-6=57", $true, $true, $false, $false, $false, $true, 1, $false, "8+23=31", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 63-6=57" }
$found = $d.Content.Find.Execute("33-26=7", $true, $true, $false, $false, $false, $true, 1, $false, "37+6=43", 2)
if ($found) { $count++ } else { Write-Output "MISSING: 33-26=7" }
Write-Output "Replaced $count of 101"